$wb = $excel.ActiveWorkbook

# --- Insert a new "Player Info" sheet as the first sheet ---
$playerInfo = $wb.Worksheets.Add($wb.Worksheets.Item(1))
$playerInfo.Name = "Player Info"

$playerInfo.Range("A1").Value = "ID"
$playerInfo.Range("B1").Value = "NAME"
$playerInfo.Range("C1").Value = "BATTING_HAND"
$playerInfo.Range("D1").Value = "BOWL_STYLE"

$headerRange = $playerInfo.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

$playerInfo.Range("A2").Value = "'4552"
$playerInfo.Range("B2").Value = "Krunal Himanshu Pandya"
$playerInfo.Range("C2").Value = "Left Handed"
$playerInfo.Range("D2").Value = "Left Arm Orthodox"

# --- Update "ODI Batting" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$batting = $wb.Worksheets.Item("ODI Batting")
$batting.Range("D1").Value = "MATCH_CODE"

$batting.Range("D2").Value = "'4454"
$batting.Range("D3").Value = "'4456"
$batting.Range("D4").Value = "'4457"
$batting.Range("D5").Value = "'4480"
$batting.Range("D6").Value = "'4482"

# --- Update "ODI Bowling" sheet: MATCH_CARD_LINK -> MATCH_CODE ---
$bowling = $wb.Worksheets.Item("ODI Bowling")
$bowling.Range("B1").Value = "MATCH_CODE"

$bowling.Range("B2").Value = "'4454"
$bowling.Range("B3").Value = "'4456"
$bowling.Range("B4").Value = "'4457"
$bowling.Range("B5").Value = "'4480"
$bowling.Range("B6").Value = "'4482"
